# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.096.66"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "1.787.63"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.17"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.53"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.284"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0703"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "2.046.86"
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("D13").Value = "1.787.61"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.80"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.623"
$ws.Range("E15").Value = "  -4.73%  "
$ws.Range("D16").Value = "34.100.34"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.16"
$ws.Range("E17").Value = "  -4.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.86"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.88"
$ws.Range("E19").Value = "  -4.55%  "
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  -5.48%  "
$ws.Range("E23").Value = "  -4.82%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.60"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.29"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("E33").Value = "  -3.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("E34").Value = "  -6.20%  "
$ws.Range("D35").Value = "1.395.89"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.70"
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.59"
$ws.Range("E42").Value = "  -5.53%  "
$ws.Range("E43").Value = "  -7.27%  "
$ws.Range("D44").Value = "0.0₆0145"
$ws.Range("E44").Value = "  +14.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.08"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.62"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("D49").Value = "1.945.46"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.18"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("E51").Value = "  -0.13%  "
